# wijzigingen mbt 2nd projectleider in zowel werkbestand als algemeen overzicht.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C for the "2e Projectleider" field,
# shifting the existing columns (Klant, Omschrijving, ...) one place to the right.
$ws.Columns("C:C").Insert() | Out-Null

# Fill in the new header cell with the shared-string label for the 2nd project leader.
$ws.Range("C2").Value = "2e Projectleider"

# Update the active selection to match the saved view.
$ws.Range("C7").Select() | Out-Null
